$wb = $excel.ActiveWorkbook

# --- Overview sheet ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Cells.Item(2, 2).Value = "Ready for handoff"       # B2: zh-cn status
$overview.Cells.Item(2, 3).Value = "Ready for handoff"       # C2: de-de status
$overview.Cells.Item(2, 4).Value = "2016-38-13 12:38:00"     # D2: Latest Handoff Date

# --- zh-cn sheet ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Cells.Item(2, 3).Value = "Ready for handoff"           # C2: Status
$zhcn.Cells.Item(2, 5).Value = "2016-03-13 12:37:57"         # E2: Latest Handoff Datetime

# --- de-de sheet ---
$dede = $wb.Worksheets.Item("de-de")
$dede.Cells.Item(2, 3).Value = "Ready for handoff"           # C2: Status
$dede.Cells.Item(2, 5).Value = "2016-03-13 12:38:00"         # E2: Latest Handoff Datetime
